$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 4

# Update row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 4

# Update row 4
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 4

# Remove rows 5-9 of data (study now uses randomization, fewer rows remain)
$ws.Range("A5:C9").ClearContents() | Out-Null

# Update selection to match final state
$ws.Range("B5:C9").Select() | Out-Null
